# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns
# with newly scraped values (also two rows whose coin/ranking swapped
# places: MultiversX <-> Algorand). A leading apostrophe is used only
# for the "Price" column because several prices are numeric-looking
# strings (e.g. "318.45", "2.693.08") that Excel would otherwise
# auto-convert to numbers; it keeps them stored as text, matching the
# original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'42.721.80"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "'2.341.19"
$ws.Range("E3").Value = "  -3.16%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'318.45"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'104.79"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  -1.86%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  -6.29%  "

$ws.Range("D10").Value = "'40.51"
$ws.Range("E10").Value = "  -3.16%  "

$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  -3.04%  "

$ws.Range("D12").Value = "'8.33"
$ws.Range("E12").Value = "  -4.21%  "

$ws.Range("D13").Value = "'0.986"
$ws.Range("E13").Value = "  -5.07%  "

$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "'15.87"
$ws.Range("E15").Value = "  -9.04%  "

$ws.Range("D16").Value = "'2.693.08"
$ws.Range("E16").Value = "  -3.13%  "

$ws.Range("D17").Value = "'2.346.03"
$ws.Range("E17").Value = "  -6.20%  "

$ws.Range("D18").Value = "'42.644.43"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  +3.86%  "

$ws.Range("E20").Value = "  -4.13%  "

$ws.Range("D21").Value = "'77.28"
$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("D22").Value = "'3.55"
$ws.Range("E22").Value = "  +1.72%  "

$ws.Range("D23").Value = "'260.41"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").Value = "'2.32"

$ws.Range("D25").Value = "'9.64"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -5.35%  "

$ws.Range("D28").Value = "'23.06"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").Value = "'174.92"
$ws.Range("E30").Value = "  -1.85%  "

$ws.Range("D31").Value = "'35.86"
$ws.Range("E31").Value = "  -6.01%  "

$ws.Range("D32").Value = "'0.0890"
$ws.Range("E32").Value = "  -5.01%  "

$ws.Range("D33").Value = "'2.99"
$ws.Range("E33").Value = "  -7.62%  "

$ws.Range("D34").Value = "'6.07"
$ws.Range("E34").Value = "  +1.56%  "

$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("E36").Value = "  +4.87%  "

$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  -7.22%  "

$ws.Range("D38").Value = "'0.0355"
$ws.Range("E38").Value = "  -4.51%  "

$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "  -5.21%  "

$ws.Range("D40").Value = "'2.63"
$ws.Range("E40").Value = "  -10.11%  "

$ws.Range("E41").Value = "  -11.94%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.232"
$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'70.37"
$ws.Range("E43").Value = "  +1.43%  "

$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").Value = "'115.18"
$ws.Range("E45").Value = "  -7.60%  "

$ws.Range("D46").Value = "'11.81"
$ws.Range("E46").Value = "  -6.72%  "

$ws.Range("D47").Value = "'5.50"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("D48").Value = "'9.17"
$ws.Range("E48").Value = "  -3.94%  "

$ws.Range("D49").Value = "'85.43"
$ws.Range("E49").Value = "  +10.84%  "

$ws.Range("D50").Value = "'73.20"
$ws.Range("E50").Value = "  +3.09%  "

$ws.Range("D51").Value = "'0.0997"
$ws.Range("E51").Value = "  -1.30%  "

